# Replace io-model/SIOM with an internal calculation (#149)
# The "SIOM" (Standard Input Output Matrix) acronym row is being removed
# from the "Key to Variables" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# Locate the row that defines the SIOM acronym (Top Level Folder = "io-model",
# Acronym = "SIOM") and delete the entire row, shifting all following rows up.
$found = $ws.Cells.Find("SIOM", [Type]::Missing, [Type]::Missing, 1)

if ($found -ne $null) {
    $rowNum = $found.Row
    $ws.Rows.Item($rowNum).Delete() | Out-Null
} else {
    # Fallback: the SIOM row is known to be row 163 in the original layout.
    $ws.Rows.Item(163).Delete() | Out-Null
}
